# Simulated Wild Card round and logged it
# Update Rushing and Receiving stat sheets with the results of the simulated game.

$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$ws = $wb.Worksheets.Item("Rushing")

# Row 2 - M.Jones
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 17

# Row 3 - D.Harris
$ws.Range("C3").Value = 130
$ws.Range("D3").Value = 65
$ws.Range("E3").Value = 17
$ws.Range("F3").Value = 48

# Row 5 - R.Stevenson
$ws.Range("C5").Value = 86
$ws.Range("D5").Value = 46
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 26

# Row 7 - B.Bolden
$ws.Range("C7").Value = 8
$ws.Range("D7").Value = 13
$ws.Range("E7").Value = 20
$ws.Range("F7").Value = 8

# Row 10 - K.Bourne
$ws.Range("C10").Value = 9
$ws.Range("D10").Value = 5

# Row 13 - J.Smith
$ws.Range("E13").Value = 1

# --- Receiving sheet ---
$ws2 = $wb.Worksheets.Item("Receiving")

# Row 2 - D.Harris
$ws2.Range("C2").Value = 22
$ws2.Range("D2").Value = 18

# Row 3 - R.Stevenson
$ws2.Range("C3").Value = 20
$ws2.Range("D3").Value = 16
$ws2.Range("G3").Value = 3
$ws2.Range("H3").Value = 2

# Row 5 - B.Bolden
$ws2.Range("C5").Value = 48
$ws2.Range("D5").Value = 39
$ws2.Range("E5").Value = 6
$ws2.Range("G5").Value = 10
$ws2.Range("H5").Value = 9

# Row 6 - N.Agholor
$ws2.Range("C6").Value = 41
$ws2.Range("D6").Value = 30
$ws2.Range("E6").Value = 25
$ws2.Range("F6").Value = 9

# Row 7 - J.Meyers
$ws2.Range("C7").Value = 107
$ws2.Range("D7").Value = 77
$ws2.Range("E7").Value = 28
$ws2.Range("F7").Value = 12
$ws2.Range("G7").Value = 17
$ws2.Range("H7").Value = 10

# Row 8 - K.Bourne
$ws2.Range("C8").Value = 62
$ws2.Range("D8").Value = 51
$ws2.Range("E8").Value = 15
$ws2.Range("F8").Value = 11
$ws2.Range("G8").Value = 9
$ws2.Range("H8").Value = 7

# Row 10 - N.Harry
$ws2.Range("E10").Value = 8

# Row 11 - Jak.Johnson
$ws2.Range("C11").Value = 5

# Row 14 - H.Henry
$ws2.Range("C14").Value = 64
$ws2.Range("D14").Value = 43
$ws2.Range("E14").Value = 16
$ws2.Range("F14").Value = 9
$ws2.Range("G14").Value = 19
